$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "65.369.55"
$ws.Cells.Item(2, 5).Value = "  -4.69%  "

$ws.Cells.Item(3, 4).Value = "3.253.92"
$ws.Cells.Item(3, 5).Value = "  -5.32%  "

$ws.Cells.Item(4, 5).Value = "  +0.14%  "

$ws.Cells.Item(5, 4).Value = "553.34"
$ws.Cells.Item(5, 5).Value = "  -3.28%  "

$ws.Cells.Item(6, 4).Value = "181.33"
$ws.Cells.Item(6, 5).Value = "  -4.25%  "

$ws.Cells.Item(7, 5).Value = "  +0.02%  "

$ws.Cells.Item(8, 4).Value = "0.589"
$ws.Cells.Item(8, 5).Value = "  -2.97%  "

$ws.Cells.Item(9, 4).Value = "3.247.54"
$ws.Cells.Item(9, 5).Value = "  -5.17%  "

$ws.Cells.Item(10, 5).Value = "  -9.01%  "

$ws.Cells.Item(11, 5).Value = "  -4.31%  "

$ws.Cells.Item(12, 4).Value = "46.98"
$ws.Cells.Item(12, 5).Value = "  -7.62%  "

$ws.Cells.Item(13, 5).Value = "  -6.65%  "

$ws.Cells.Item(14, 4).Value = "628.17"
$ws.Cells.Item(14, 5).Value = "  +0.01%  "

$ws.Cells.Item(15, 4).Value = "8.51"
$ws.Cells.Item(15, 5).Value = "  -5.51%  "

$ws.Cells.Item(16, 4).Value = "3.781.65"
$ws.Cells.Item(16, 5).Value = "  -4.89%  "

$ws.Cells.Item(17, 4).Value = "65.307.81"
$ws.Cells.Item(17, 5).Value = "  -4.44%  "

$ws.Cells.Item(18, 5).Value = "  -3.27%  "

$ws.Cells.Item(19, 4).Value = "17.67"
$ws.Cells.Item(19, 5).Value = "  -1.80%  "

$ws.Cells.Item(20, 4).Value = "3.257.05"
$ws.Cells.Item(20, 5).Value = "  -4.80%  "

$ws.Cells.Item(21, 5).Value = "  -6.93%  "

$ws.Cells.Item(22, 4).Value = "0.898"
$ws.Cells.Item(22, 5).Value = "  -3.82%  "

$ws.Cells.Item(23, 4).Value = "17.71"
$ws.Cells.Item(23, 5).Value = "  +0.06%  "

$ws.Cells.Item(24, 4).Value = "105.75"
$ws.Cells.Item(24, 5).Value = "  +7.66%  "

$ws.Cells.Item(25, 4).Value = "4.90"
$ws.Cells.Item(25, 5).Value = "  -7.86%  "

$ws.Cells.Item(26, 4).Value = "3.94"
$ws.Cells.Item(26, 5).Value = "  -6.55%  "

$ws.Cells.Item(27, 5).Value = "  -5.98%  "

$ws.Cells.Item(28, 4).Value = "9.47"
$ws.Cells.Item(28, 5).Value = "  -3.17%  "

$ws.Cells.Item(29, 4).Value = "8.62"
$ws.Cells.Item(29, 5).Value = "  -5.66%  "

$ws.Cells.Item(30, 4).Value = "30.20"
$ws.Cells.Item(30, 5).Value = "  -5.88%  "

$ws.Cells.Item(31, 4).Value = "3.93"
$ws.Cells.Item(31, 5).Value = "  -4.86%  "

$ws.Cells.Item(32, 4).Value = "6.23"
$ws.Cells.Item(32, 5).Value = "  -6.06%  "

$ws.Cells.Item(33, 4).Value = "10.96"
$ws.Cells.Item(33, 5).Value = "  -4.56%  "

$ws.Cells.Item(34, 4).Value = "543.82"
$ws.Cells.Item(34, 5).Value = "  +10.08%  "

$ws.Cells.Item(35, 5).Value = "  -3.36%  "

$ws.Cells.Item(36, 4).Value = "0.999"
$ws.Cells.Item(36, 5).Value = "  -0.05%  "

$ws.Cells.Item(37, 4).Value = "56.95"
$ws.Cells.Item(37, 5).Value = "  -5.80%  "

$ws.Cells.Item(38, 4).Value = "3.613.79"

$ws.Cells.Item(39, 4).Value = "3.37"
$ws.Cells.Item(39, 5).Value = "  -2.39%  "

$ws.Cells.Item(40, 2).Value = "Kaspa"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(40, 4).Value = "0.129"
$ws.Cells.Item(40, 5).Value = "  -1.79%  "

$ws.Cells.Item(41, 5).Value = "  -7.59%  "

$ws.Cells.Item(42, 2).Value = "Fetch.AI"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(42, 4).Value = "2.71"
$ws.Cells.Item(42, 5).Value = "  -5.44%  "

$ws.Cells.Item(43, 2).Value = "CoreDAO"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Cells.Item(43, 4).Value = "3.34"
$ws.Cells.Item(43, 5).Value = "  -7.08%  "

$ws.Cells.Item(44, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(44, 4).Value = "32.08"
$ws.Cells.Item(44, 5).Value = "  -5.55%  "

$ws.Cells.Item(45, 4).Value = "0.332"
$ws.Cells.Item(45, 5).Value = "  -9.16%  "

$ws.Cells.Item(46, 4).Value = "3.28"
$ws.Cells.Item(46, 5).Value = "  -1.45%  "

$ws.Cells.Item(47, 4).Value = "0.0413"
$ws.Cells.Item(47, 5).Value = "  -4.93%  "

$ws.Cells.Item(48, 5).Value = "  -6.63%  "

$ws.Cells.Item(49, 5).Value = "  -3.36%  "

$ws.Cells.Item(50, 4).Value = "0.999"
$ws.Cells.Item(50, 5).Value = "  +0.11%  "

$ws.Cells.Item(51, 5).Value = "  +1.23%  "
